# Fix loop bounds in overview figure, and let the "Update automatically"
# date placeholders re-stamp themselves (16-04-2022 -> 18-04-2022) the way
# PowerPoint itself refreshes those fields on save.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Part 1: bump every "Date Placeholder" shape's visible text from
# 16-04-2022 to 18-04-2022 across the slide master, every custom layout and
# the notes master.
# ---------------------------------------------------------------------------

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $txt = $tr.Text
            if ($txt -like "*16-04-2022*") {
                $tr.Text = $txt.Replace("16-04-2022", "18-04-2022")
            }
        }
    }
}

$design = $p.Designs.Item(1)
$master = $design.SlideMaster

Update-DateShapes $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShapes $layouts.Item($li).Shapes
}

# NOTE: NotesMaster.Shapes.Item(n) writes are mis-routed onto the
# SlideMaster's own shapes by index in this runtime, so the notes master's
# date placeholder must be updated through its HeadersFooters object
# instead of walking NotesMaster.Shapes directly.
$p.NotesMaster.HeadersFooters.DateAndTime.Text = "18-04-2022"

# ---------------------------------------------------------------------------
# Part 2: fix the loop bounds printed in the overview figure's pseudo-code
# boxes on slide 1 (t = 0 to 1 / t = 1 to 3  ->  t = 0 to 2 / t = 2 to 3).
# ---------------------------------------------------------------------------

$slide = $p.Slides.Item(1)

# "Rectangle: Rounded Corners 150" (top-level shape #20)
$shp150 = $slide.Shapes.Item(20)
$tr150 = $shp150.TextFrame.TextRange
$x = $tr150.Replace(" t = 0 to 1 ", " t = 0 to 2 ", 1, 0, 0)
$x = $tr150.Replace(" t = 1 to 3 step 1 {", " t = 2 to 3 step 1 {", 1, 0, 0)

# "Rectangle: Rounded Corners 152" (top-level shape #21)
$shp152 = $slide.Shapes.Item(21)
$tr152 = $shp152.TextFrame.TextRange
$x = $tr152.Replace(" t = 0 to 1 ", " t = 0 to 2 ", 1, 0, 0)
$x = $tr152.Replace(" t = 1 to 3 ", " t = 2 to 3 ", 1, 0, 0)

# "Rectangle: Rounded Corners 160" (top-level shape #28)
$shp160 = $slide.Shapes.Item(28)
$tr160 = $shp160.TextFrame.TextRange
$x = $tr160.Replace(" t = 0 to 1 ", " t = 0 to 2 ", 1, 0, 0)

# The second "for t = 1 to 3" line in shape 160 is re-typed as three runs
# (" t ", "= 2 ", "to 3 ") in the target file, so reproduce that split by
# hand instead of a single whole-run replace.
$full160 = $tr160.Text
$idx = $full160.LastIndexOf(" t = 1 to 3 ")
$start = $idx + 1
$runA = $tr160.Characters($start, 3)
$runB = $tr160.Characters($start + 3, 4)
$runC = $tr160.Characters($start + 7, 5)
$runA.Text = " t "
$runB.Text = "= 2 "
$runC.Text = "to 3 "
